# Add a new "alcohol_content" (ABV) column as column D to the beverage
# report sheet, matching the price column's header style for the header
# cell, and leaving the non-alcoholic drinks (rows 19-23) with a present
# but blank value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell (D1) -------------------------------------------------
$ws.Range("D1").Value = "alcohol_content"
# Match the look of the existing header row (bold, bordered, centered)
# by copying the formatting of the neighboring header cell.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats

# --- Alcohol content values for each beverage (rows 2-18) -------------
$abv = @{
    2  = 13.5
    3  = 12
    4  = 13
    5  = 12.5
    6  = 14
    7  = 11
    8  = 13.5
    9  = 13.8
    10 = 14.2
    11 = 12
    12 = 5
    13 = 8
    14 = 5
    15 = 4.5
    16 = 8
    17 = 6.8
    18 = 5.9
}

foreach ($row in $abv.Keys) {
    $ws.Cells.Item($row, 4).Value = $abv[$row]
}

# --- Non-alcoholic beverages (rows 19-23): present but blank value ----
19..23 | ForEach-Object {
    $r = $_
    $ws.Cells.Item($r, 4).Value = "'"
    $ws.Range("C" + $r).Copy()
    $ws.Range("D" + $r).PasteSpecial(-4122)   # xlPasteFormats (keep unstyled like column C)
}

Write-Output "alcohol_content column added"
